$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 17 sim play logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "8 2 6 9 1 8 -3 1 3 1 3 3 3 7 3 2 2 11 2 6 15 25 2 3 46 4 6 1 10 1 -1 2 5 4 0 8 4 4 4 1 4 15 2 5 0 0 7 3 4 2 8 1 0 4 6 8 1 -2 4 4 8 2 3 4 3 2 2 0 3 0 7 6 5 1 6 9 1 2 2 3 4 3 1 8 0 17 8 -5 7 9 6 3 0 8 5 41 12 7 7 1 5 1 4 3 8 2 5 -1 4 0 8 22 2 11 9 9 8 2 -1 3 5 -1 0 0 3 1 4 12 4 6 6 0 8 3 1 3 14 3 7 3 -3 7 5 2 1 7 5 2 3 4 0 4 5 0 4 5 4 -1 6 34 1 0 -1 -2 6 5 7 13 2 6 2 12 3 -1 4 0 7 15 12 -3 1 22 4 -1 1 7 5 4 2 3 -1 15 0 -4 4 0 4 11 7 13 24 8 15 0 1 7 6 3 9 2 3 10 0 3 8 9 4 28 3 8 5 2 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1 14 -3 10 22 1 29 18 4 11 5 18 23 2 0 7 8 3 8 2 16 6 3 0 2 2 26 3 4 1 -2 12 3 1 2 4 1 8 7 4 5 2 1 3 3 5 1 5 3 25 2 6 2 -1 8 3 7 3 2 11 1 2 2 7 2 0 0 5"
$ws.Range("C2").Value = "2 4 2 2 -4 2 0 9 2 1 2 25 -2 5 18 3 8 -1 6 18 -1 2 -1 6 2 1 3 2 8 3 4 -2 6 1 12 1 8 2 4 2 7 -1 4 1 -1 11 4 -4 7 3 1 7 4 7 6 0 1 0 5 1 10 3 4 0 12 0 -6 2 -1 2 4 1 3 2 4 2 5 4 2 10 5 2 3 15 3 6 4 1 9 5 0 -1 23 9 -1 -2 9 5 9 0 -1 2 4 3 76 2 -3 7 4 6 3 19 0 4 3 1 2 3 0 -1 13 1 7 5 5 5 -3 2 2 5 3 -2 12 3 1 1 0 2 6 2 0 3 1 7 6 3 2 3 -4 -1 15 1 -2 3 0 4 11 7 -1 5 1 6 2 8 2 0 -2 3 5 5 15 6 5 4 -4 5 -2 -2 5 8 5 13 0 0 0 5 1 -1 3 8 10 6 10 3 3 11 2 4 5 1 6 0 18 4 3 7 -1 2 2 3 5 40 10 9 0 3 2 1 19 6 5 18 9 3 -3 1 5 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2 0 3 47 4 2 12 2 3 1 5 -1 8 2 1 1 2 9 1 -2 2 0 7 21 1 4 3 1 1 11 4 21 3 -2 0 7 8 -4 5 2 2 0 2 30 3 4 1 13 -1 1 7 1 7 8 5 7 5 0 9 5 9 0 0 2 1 16 1 7 31 2 0 2 2 1 3 1 9 6 8 11 3"
$ws.Range("B3").Value = "7 -2 6 5 10 16 13 4 7 12 10 5 37 10 8 3 7 9 9 12 7 4 6 12 -1 13 10 13 5 1 35 5 6 1 7 7 22 13 6 7 8 6 9 2 41 3 23 11 13 1 28 3 7 9 7 17 4 7 14 18 2 41 7 13 22 11 5 11 9 5 6 13 10 11 6 5 19 7 5 25 31 9 20 1 16 7 8 5 16 21 2 37 9 5 5 1 19 23 8 24 35 61 15 53 41 1 -2 11 5 16 16 8 9 12 12 6 4 3 11 6 7 3 3 21 19 2 14 5 17 12 29 1 5 9 11 18 9 15 3 4 1 6 6 31 12 7 20 1 11 8 1 20 3 15 7 8 6 15 4 10 6 3 3 8 8 17 11 14 19 13 1 11 -1 20 8 4 3 7 1 16 0 4 9 9 14 1 7 2 9 11 13 7 9 5 22 8 12 15 3 5 -2 28 1 7 19 5 14 3 11 17 15 4 7 7 12 3 6 8 6 6 16 57 12 23 31 49 25 8 36 20 43 7 11 9 7 20 2 18 3 6 14 12 7 6 10 6 7 16 7 31 1 9 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6 1 10 11 5 3 25 24 21 9 5 6 10 12 6 1 13 9 -2 10 2 4 10 15 9 12 12 3 4 4 4 2 15 1 24 3 12 9 15 9 10 12 11 4 10 9 7 11 7 7 7 8 28 20 11 8 7 13 7 8 12 9 4 13 3 10 1 12 9 11 18 9 11 11 6 5 28 9 17 15 19 2"
$ws.Range("C3").Value = "14 6 6 12 9 10 6 19 22 5 9 24 4 -2 5 24 14 1 13 7 6 9 5 18 12 4 5 14 4 11 16 11 3 5 8 0 4 2 7 8 4 5 6 0 11 17 10 73 14 7 4 15 1 11 15 37 3 8 7 8 5 8 12 3 6 12 12 6 15 11 9 7 8 5 8 10 6 26 8 12 8 11 5 4 3 11 8 17 6 8 5 5 4 1 8 6 4 4 11 11 5 48 7 8 8 8 4 13 20 6 6 14 22 5 8 9 23 7 8 0 14 15 22 12 8 3 8 5 1 -1 6 12 10 4 10 -2 11 4 6 12 -1 6 7 5 5 4 5 28 9 4 16 9 4 6 7 3 3 2 9 10 4 21 18 14 28 19 14 12 11 11 6 9 -7 21 18 8 -4 18 18 5 23 9 18 0 3 8 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 12 7 6 15 5 7 5 9 20 8 16 23 18 5 13 16 13 7 5 20 5 6 24 4 18 0 8 5 11 10 6 -7 14 58 1 9 8 8 9 4 11 13 9 16 4 23 5 6 8 5 8 9 13 5 14 6 10 9 3 6 9 13 21"

# --- OFF sheet ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 143
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = 51
$ws.Range("J2").Value = 29
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 10
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 224
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 108
$ws.Range("G3").Value = 28
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 60
$ws.Range("L3").Value = 313
$ws.Range("M3").Value = 210
$ws.Range("Q3").Value = 518

# --- DEF sheet ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 177
$ws.Range("D2").Value = 8
$ws.Range("F2").Value = 49
$ws.Range("G2").Value = 53
$ws.Range("J2").Value = 30
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 31
$ws.Range("P2").Value = 15
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 143
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 89
$ws.Range("H3").Value = 33
$ws.Range("I3").Value = 55
$ws.Range("J3").Value = 43
$ws.Range("L3").Value = 308
$ws.Range("M3").Value = 179
$ws.Range("Q3").Value = 509

# --- ST sheet ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 86
$ws.Range("F2").Value = 104
$ws.Range("G2").Value = 102
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 41
$ws.Range("K2").Value = 37
$ws.Range("B3").Value = 47
$ws.Range("B4").Value = "70 65 63 56 63 62 58 66 68 70 63 70 54 61 58 64 63 63 61 62 40 63 66 60 66 63 64 62 63 62 65 56 61 56 58 57 58 54 63"
$ws.Range("B5").Value = "31 21 0 19 17 24 16 16 31 33 0 28 12 17 17 20 17 21 31 18 7 33 18 29 33 33 19 14 15 18 14 15 33 7 21 10 18 21 17"
$ws.Range("B6").Value = "75 23 19 26 0 20 0 30 18 19 18 45 33 28 23 22 12 12 27 20 17 17 27 29 13 8 28 18 31"
$ws.Range("D5").Value = "0 1 0 0 0 20 0 0 0 13 6 33 0 0 5 0 9 6 12 -1 7 0 6 14 16 0 9 0 0 0 0 0 0 0 26 0 0 11 10 1 18 8 0 0 0 0 0 4 0 0 0 5 17 14 0 0 4 15 15"

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 12

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B3").Value = 23
$ws.Range("D3").Value = 11
$ws.Range("B4").Value = 1
